# Updated cryptos list data refresh (prices & volume %) for Sheet1.
# Source data is plain text (inline strings); some new price values parse as
# numbers (e.g. "328.20", "2.22"), so force those specific cells to Text format
# BEFORE assigning the value, to avoid Excel auto-converting them to numbers
# (which would lose trailing zeros / exact formatting of the original text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new value would otherwise be
# auto-detected as a number by Excel.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Assign the updated cell values (coin name/link swap for rows 30-31, plus
# refreshed Price (D) and Volume(1h) (E) columns for all coins).
$ws.Range("D2").Value = "43.617.79"
$ws.Range("E2").Value = "  +2.96%  "
$ws.Range("D3").Value = "2.418.44"
$ws.Range("E3").Value = "  +8.80%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "328.20"
$ws.Range("E5").Value = "  +13.33%  "
$ws.Range("D6").Value = "105.10"
$ws.Range("E6").Value = "  -5.30%  "
$ws.Range("E7").Value = "  +4.05%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "0.671"
$ws.Range("E9").Value = "  +12.32%  "
$ws.Range("D10").Value = "42.45"
$ws.Range("E10").Value = "  -2.34%  "
$ws.Range("D11").Value = "0.0948"
$ws.Range("E11").Value = "  +4.22%  "
$ws.Range("D12").Value = "8.69"
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("E13").Value = "  +3.25%  "
$ws.Range("D14").Value = "17.34"
$ws.Range("E14").Value = "  +16.60%  "
$ws.Range("E15").Value = "  +2.84%  "
$ws.Range("D16").Value = "2.788.45"
$ws.Range("E16").Value = "  +8.99%  "
$ws.Range("D17").Value = "2.419.91"
$ws.Range("E17").Value = "  +8.45%  "
$ws.Range("D18").Value = "43.676.89"
$ws.Range("E18").Value = "  +3.20%  "
$ws.Range("D19").Value = "7.55"
$ws.Range("E19").Value = "  +6.31%  "
$ws.Range("D20").Value = "0.0000110"
$ws.Range("E20").Value = "  +5.53%  "
$ws.Range("D21").Value = "76.20"
$ws.Range("E21").Value = "  +4.85%  "
$ws.Range("E22").Value = "  +5.85%  "
$ws.Range("D23").Value = "274.15"
$ws.Range("E23").Value = "  +19.06%  "
$ws.Range("D24").Value = "2.47"
$ws.Range("E24").Value = "  +2.51%  "
$ws.Range("D25").Value = "9.65"
$ws.Range("E25").Value = "  +7.77%  "
$ws.Range("D26").Value = "12.07"
$ws.Range("E26").Value = "  +5.93%  "
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("D29").Value = "23.13"
$ws.Range("E29").Value = "  +11.06%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "178.36"
$ws.Range("E30").Value = "  +2.91%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "2.22"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D32").Value = "38.01"
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("E33").Value = "  +3.96%  "
$ws.Range("D34").Value = "0.0941"
$ws.Range("E34").Value = "  +7.09%  "
$ws.Range("D35").Value = "5.99"
$ws.Range("E35").Value = "  +6.97%  "
$ws.Range("E36").Value = "  +6.70%  "
$ws.Range("D37").Value = "4.90"
$ws.Range("E37").Value = "  -0.81%  "
$ws.Range("D38").Value = "0.0373"
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("D39").Value = "4.07"
$ws.Range("E39").Value = "  -3.73%  "
$ws.Range("D40").Value = "0.110"
$ws.Range("E40").Value = "  +5.49%  "
$ws.Range("E41").Value = "  +21.75%  "
$ws.Range("E42").Value = "  +25.53%  "
$ws.Range("D43").Value = "129.24"
$ws.Range("E43").Value = "  +27.58%  "
$ws.Range("D44").Value = "0.237"
$ws.Range("E44").Value = "  +2.83%  "
$ws.Range("D45").Value = "70.36"
$ws.Range("E45").Value = "  -4.23%  "
$ws.Range("D46").Value = "12.81"
$ws.Range("E46").Value = "  +4.55%  "
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("D48").Value = "9.78"
$ws.Range("E48").Value = "  +16.44%  "
$ws.Range("D49").Value = "5.73"
$ws.Range("E49").Value = "  +7.65%  "
$ws.Range("D50").Value = "87.38"
$ws.Range("E50").Value = "  +67.86%  "
$ws.Range("D51").Value = "1.33"
$ws.Range("E51").Value = "  +4.57%  "

